# Updated viz for Day 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Day 6 row (row 10): title + times
$ws.Range("B10").Value = "Day 6: Universal Orbit Map"
$ws.Range("C10").Value = 0.0067245370370370367
$ws.Range("E10").Value = 0.013842592592592594
$ws.Range("F10").Value = 0.013842592592592594
$ws.Range("H10").Value = "1st"

# The "Adj. Median" row is no longer used now that every day has data
$ws.Range("D33:E33").ClearContents()

# Update selection to match the author's final cursor position
$ws.Range("B34").Select() | Out-Null
